$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) and Volume (E) columns keep their original text
# formatting (leading zeros, trailing zeros, thousand-dot separators,
# scientific-notation-prone small decimals, etc.) by forcing the cells to
# Text format before writing the literal string value into them.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.832.22'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.814.65'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -1.20%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.46%  '
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -1.95%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4598'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -2.95%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3630'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -1.68%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07208'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -3.38%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.8566'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -3.26%  '
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -3.74%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07507'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +2.27%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.747.25'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -7.15%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.311'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -2.50%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.508'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -1.09%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '91.44'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -1.92%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.008'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -0.28%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008560'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -2.94%  '
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '26.879.22'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -2.13%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.38'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -2.93%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.127'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -3.58%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.47'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -2.05%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.956.27'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -7.10%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '150.99'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -0.74%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.847'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -2.89%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.09'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -2.97%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.066'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -3.61%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.071'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -3.52%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '115.07'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -2.23%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08844'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -1.73%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.946'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.400'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -3.32%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.126'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -4.50%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.7119'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -5.97%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.003'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.92%  '
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -2.99%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.05237'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -1.99%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.410'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +0.30%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01913'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -2.20%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.919'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -2.41%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '7.127'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -2.60%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.5122'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -3.93%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.1615'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -2.77%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.137'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -4.43%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.4776'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -2.67%  '
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -0.47%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '102.81'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -2.07%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '10.01'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -5.15%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06231'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -1.20%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.611'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -3.94%  '
